$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 22.95531766666667
$ws.Range("H2").Value = 68.865953
$ws.Range("I2").Value = 0.1720020945576478
$ws.Range("J2").Value = 0.1720020945576478
$ws.Range("M2").Value = 159.4836373333333
$ws.Range("N2").Value = 478.450912
$ws.Range("O2").Value = 0.2983285084902258
$ws.Range("P2").Value = 0.2983285084902258
$ws.Range("Q2").Value = 3660.997557622127
$ws.Range("R2").Value = 32948.97801859914
$ws.Range("S2").Value = 0.05131312832657786
$ws.Range("T2").Value = 0.05131312832657786
$ws.Range("G3").Value = 22.95531766666667
$ws.Range("H3").Value = 68.865953
$ws.Range("I3").Value = 0.1720020945576478
$ws.Range("J3").Value = 0.1720020945576478
$ws.Range("O3").Value = 0.3227862111630279
$ws.Range("P3").Value = 0.3227862111630279
$ws.Range("Q3").Value = 3961.135114717547
$ws.Range("R3").Value = 35650.21603245792
$ws.Range("S3").Value = 0.05551990441436801
$ws.Range("T3").Value = 0.05551990441436801
$ws.Range("G4").Value = 22.95531766666667
$ws.Range("H4").Value = 68.865953
$ws.Range("I4").Value = 0.1720020945576478
$ws.Range("J4").Value = 0.1720020945576478
$ws.Range("M4").Value = 74.38770566666666
$ws.Range("N4").Value = 223.163117
$ws.Range("O4").Value = 0.1391489036280481
$ws.Range("P4").Value = 0.1391489036280482
$ws.Range("Q4").Value = 1707.593414072834
$ws.Range("R4").Value = 15368.3407266555
$ws.Range("S4").Value = 0.02393390287942456
$ws.Range("T4").Value = 0.02393390287942456
$ws.Range("G5").Value = 22.95531766666667
$ws.Range("H5").Value = 68.865953
$ws.Range("I5").Value = 0.1720020945576478
$ws.Range("J5").Value = 0.1720020945576478
$ws.Range("M5").Value = 58.41461433333333
$ws.Range("N5").Value = 175.243843
$ws.Range("O5").Value = 0.1092697975759847
$ws.Range("P5").Value = 0.1092697975759848
$ws.Range("Q5").Value = 1340.926028397487
$ws.Range("R5").Value = 12068.33425557738
$ws.Range("S5").Value = 0.01879463405495957
$ws.Range("T5").Value = 0.01879463405495957
$ws.Range("G6").Value = 22.95531766666667
$ws.Range("H6").Value = 68.865953
$ws.Range("I6").Value = 0.1720020945576478
$ws.Range("J6").Value = 0.1720020945576478
$ws.Range("M6").Value = 69.746216
$ws.Range("N6").Value = 209.238648
$ws.Range("O6").Value = 0.1304665791427133
$ws.Range("P6").Value = 0.1304665791427133
$ws.Range("Q6").Value = 1601.04654432795
$ws.Range("R6").Value = 14409.41889895155
$ws.Range("S6").Value = 0.02244052488231782
$ws.Range("T6").Value = 0.02244052488231783
$ws.Range("I7").Value = 0.4661646602805707
$ws.Range("J7").Value = 0.4661646602805707
$ws.Range("M7").Value = 159.4836373333333
$ws.Range("N7").Value = 478.450912
$ws.Range("O7").Value = 0.2983285084902258
$ws.Range("P7").Value = 0.2983285084902258
$ws.Range("Q7").Value = 9922.133141029444
$ws.Range("R7").Value = 89299.19826926499
$ws.Range("S7").Value = 0.1390702078123555
$ws.Range("T7").Value = 0.1390702078123555
$ws.Range("I8").Value = 0.4661646602805707
$ws.Range("J8").Value = 0.4661646602805707
$ws.Range("O8").Value = 0.3227862111630279
$ws.Range("P8").Value = 0.3227862111630279
$ws.Range("S8").Value = 0.1504715244700655
$ws.Range("T8").Value = 0.1504715244700655
$ws.Range("I9").Value = 0.4661646602805707
$ws.Range("J9").Value = 0.4661646602805707
$ws.Range("M9").Value = 74.38770566666666
$ws.Range("N9").Value = 223.163117
$ws.Range("O9").Value = 0.1391489036280481
$ws.Range("P9").Value = 0.1391489036280482
$ws.Range("Q9").Value = 4627.965175748544
$ws.Range("R9").Value = 41651.6865817369
$ws.Range("S9").Value = 0.06486630138818293
$ws.Range("T9").Value = 0.06486630138818295
$ws.Range("I10").Value = 0.4661646602805707
$ws.Range("J10").Value = 0.4661646602805707
$ws.Range("M10").Value = 58.41461433333333
$ws.Range("N10").Value = 175.243843
$ws.Range("O10").Value = 0.1092697975759847
$ws.Range("P10").Value = 0.1092697975759848
$ws.Range("Q10").Value = 3634.213455928495
$ws.Range("R10").Value = 32707.92110335646
$ws.Range("S10").Value = 0.05093771806593566
$ws.Range("T10").Value = 0.05093771806593566
$ws.Range("I11").Value = 0.4661646602805707
$ws.Range("J11").Value = 0.4661646602805707
$ws.Range("M11").Value = 69.746216
$ws.Range("N11").Value = 209.238648
$ws.Range("O11").Value = 0.1304665791427133
$ws.Range("P11").Value = 0.1304665791427133
$ws.Range("Q11").Value = 4339.199010043885
$ws.Range("R11").Value = 39052.79109039497
$ws.Range("S11").Value = 0.06081890854403114
$ws.Range("T11").Value = 0.06081890854403115
$ws.Range("G12").Value = 7.783044333333334
$ws.Range("H12").Value = 23.349133
$ws.Range("I12").Value = 0.05831763893698088
$ws.Range("J12").Value = 0.05831763893698089
$ws.Range("M12").Value = 159.4836373333333
$ws.Range("N12").Value = 478.450912
$ws.Range("O12").Value = 0.2983285084902258
$ws.Range("P12").Value = 0.2983285084902258
$ws.Range("Q12").Value = 1241.268219806589
$ws.Range("R12").Value = 11171.4139782593
$ws.Range("S12").Value = 0.01739781424274102
$ws.Range("T12").Value = 0.01739781424274103
$ws.Range("G13").Value = 7.783044333333334
$ws.Range("H13").Value = 23.349133
$ws.Range("I13").Value = 0.05831763893698088
$ws.Range("J13").Value = 0.05831763893698089
$ws.Range("O13").Value = 0.3227862111630279
$ws.Range("P13").Value = 0.3227862111630279
$ws.Range("Q13").Value = 1343.030432244367
$ws.Range("R13").Value = 12087.2738901993
$ws.Range("S13").Value = 0.01882412971644153
$ws.Range("T13").Value = 0.01882412971644153
$ws.Range("G14").Value = 7.783044333333334
$ws.Range("H14").Value = 23.349133
$ws.Range("I14").Value = 0.05831763893698088
$ws.Range("J14").Value = 0.05831763893698089
$ws.Range("M14").Value = 74.38770566666666
$ws.Range("N14").Value = 223.163117
$ws.Range("O14").Value = 0.1391489036280481
$ws.Range("P14").Value = 0.1391489036280482
$ws.Range("Q14").Value = 578.9628110586178
$ws.Range("R14").Value = 5210.665299527562
$ws.Range("S14").Value = 0.00811483552025726
$ws.Range("T14").Value = 0.008114835520257262
$ws.Range("G15").Value = 7.783044333333334
$ws.Range("H15").Value = 23.349133
$ws.Range("I15").Value = 0.05831763893698088
$ws.Range("J15").Value = 0.05831763893698089
$ws.Range("M15").Value = 58.41461433333333
$ws.Range("N15").Value = 175.243843
$ws.Range("O15").Value = 0.1092697975759847
$ws.Range("P15").Value = 0.1092697975759848
$ws.Range("Q15").Value = 454.6435330709022
$ws.Range("R15").Value = 4091.791797638119
$ws.Range("S15").Value = 0.006372356601753267
$ws.Range("T15").Value = 0.006372356601753269
$ws.Range("G16").Value = 7.783044333333334
$ws.Range("H16").Value = 23.349133
$ws.Range("I16").Value = 0.05831763893698088
$ws.Range("J16").Value = 0.05831763893698089
$ws.Range("M16").Value = 69.746216
$ws.Range("N16").Value = 209.238648
$ws.Range("O16").Value = 0.1304665791427133
$ws.Range("P16").Value = 0.1304665791427133
$ws.Range("Q16").Value = 542.8378912102427
$ws.Range("R16").Value = 4885.541020892185
$ws.Range("S16").Value = 0.007608502855787795
$ws.Range("T16").Value = 0.007608502855787798
$ws.Range("G17").Value = 30.44016466666666
$ws.Range("H17").Value = 91.320494
$ws.Range("I17").Value = 0.2280853681650076
$ws.Range("J17").Value = 0.2280853681650076
$ws.Range("M17").Value = 159.4836373333333
$ws.Range("N17").Value = 478.450912
$ws.Range("O17").Value = 0.2983285084902258
$ws.Range("P17").Value = 0.2983285084902258
$ws.Range("Q17").Value = 4854.708182065614
$ws.Range("R17").Value = 43692.37363859053
$ws.Range("S17").Value = 0.06804436769311074
$ws.Range("T17").Value = 0.06804436769311074
$ws.Range("G18").Value = 30.44016466666666
$ws.Range("H18").Value = 91.320494
$ws.Range("I18").Value = 0.2280853681650076
$ws.Range("J18").Value = 0.2280853681650076
$ws.Range("O18").Value = 0.3227862111630279
$ws.Range("P18").Value = 0.3227862111630279
$ws.Range("Q18").Value = 5252.709063312505
$ws.Range("R18").Value = 47274.38156981255
$ws.Range("S18").Value = 0.07362281181170711
$ws.Range("T18").Value = 0.07362281181170711
$ws.Range("G19").Value = 30.44016466666666
$ws.Range("H19").Value = 91.320494
$ws.Range("I19").Value = 0.2280853681650076
$ws.Range("J19").Value = 0.2280853681650076
$ws.Range("M19").Value = 74.38770566666666
$ws.Range("N19").Value = 223.163117
$ws.Range("O19").Value = 0.1391489036280481
$ws.Range("P19").Value = 0.1391489036280482
$ws.Range("Q19").Value = 2264.374009668866
$ws.Range("R19").Value = 20379.3660870198
$ws.Range("S19").Value = 0.03173782891376052
$ws.Range("T19").Value = 0.03173782891376053
$ws.Range("G20").Value = 30.44016466666666
$ws.Range("H20").Value = 91.320494
$ws.Range("I20").Value = 0.2280853681650076
$ws.Range("J20").Value = 0.2280853681650076
$ws.Range("M20").Value = 58.41461433333333
$ws.Range("N20").Value = 175.243843
$ws.Range("O20").Value = 0.1092697975759847
$ws.Range("P20").Value = 0.1092697975759848
$ws.Range("Q20").Value = 1778.150479246493
$ws.Range("R20").Value = 16003.35431321844
$ws.Range("S20").Value = 0.02492284200943434
$ws.Range("T20").Value = 0.02492284200943434
$ws.Range("G21").Value = 30.44016466666666
$ws.Range("H21").Value = 91.320494
$ws.Range("I21").Value = 0.2280853681650076
$ws.Range("J21").Value = 0.2280853681650076
$ws.Range("M21").Value = 69.746216
$ws.Range("N21").Value = 209.238648
$ws.Range("O21").Value = 0.1304665791427133
$ws.Range("P21").Value = 0.1304665791427133
$ws.Range("Q21").Value = 2123.086299916901
$ws.Range("R21").Value = 19107.77669925211
$ws.Range("S21").Value = 0.02975751773699487
$ws.Range("T21").Value = 0.02975751773699488
$ws.Range("G22").Value = 10.06688366666667
$ws.Range("H22").Value = 30.200651
$ws.Range("I22").Value = 0.07543023805979308
$ws.Range("J22").Value = 0.07543023805979308
$ws.Range("M22").Value = 159.4836373333333
$ws.Range("N22").Value = 478.450912
$ws.Range("O22").Value = 0.2983285084902258
$ws.Range("P22").Value = 0.2983285084902258
$ws.Range("Q22").Value = 1605.503223771524
$ws.Range("R22").Value = 14449.52901394371
$ws.Range("S22").Value = 0.02250299041544073
$ws.Range("T22").Value = 0.02250299041544073
$ws.Range("G23").Value = 10.06688366666667
$ws.Range("H23").Value = 30.200651
$ws.Range("I23").Value = 0.07543023805979308
$ws.Range("J23").Value = 0.07543023805979308
$ws.Range("O23").Value = 0.3227862111630279
$ws.Range("P23").Value = 0.3227862111630279
$ws.Range("Q23").Value = 1737.126314993849
$ws.Range("R23").Value = 15634.13683494464
$ws.Range("S23").Value = 0.02434784075044584
$ws.Range("T23").Value = 0.02434784075044584
$ws.Range("G24").Value = 10.06688366666667
$ws.Range("H24").Value = 30.200651
$ws.Range("I24").Value = 0.07543023805979308
$ws.Range("J24").Value = 0.07543023805979308
$ws.Range("M24").Value = 74.38770566666666
$ws.Range("N24").Value = 223.163117
$ws.Range("O24").Value = 0.1391489036280481
$ws.Range("P24").Value = 0.1391489036280482
$ws.Range("Q24").Value = 748.8523791765741
$ws.Range("R24").Value = 6739.671412589167
$ws.Range("S24").Value = 0.01049603492642288
$ws.Range("T24").Value = 0.01049603492642288
$ws.Range("G25").Value = 10.06688366666667
$ws.Range("H25").Value = 30.200651
$ws.Range("I25").Value = 0.07543023805979308
$ws.Range("J25").Value = 0.07543023805979308
$ws.Range("M25").Value = 58.41461433333333
$ws.Range("N25").Value = 175.243843
$ws.Range("O25").Value = 0.1092697975759847
$ws.Range("P25").Value = 0.1092697975759848
$ws.Range("Q25").Value = 588.053126926866
$ws.Range("R25").Value = 5292.478142341794
$ws.Range("S25").Value = 0.008242246843901931
$ws.Range("T25").Value = 0.008242246843901933
$ws.Range("G26").Value = 10.06688366666667
$ws.Range("H26").Value = 30.200651
$ws.Range("I26").Value = 0.07543023805979308
$ws.Range("J26").Value = 0.07543023805979308
$ws.Range("M26").Value = 69.746216
$ws.Range("N26").Value = 209.238648
$ws.Range("O26").Value = 0.1304665791427133
$ws.Range("P26").Value = 0.1304665791427133
$ws.Range("Q26").Value = 702.1270426622054
$ws.Range("R26").Value = 6319.143383959849
$ws.Range("S26").Value = 0.0098411251235817
$ws.Range("T26").Value = 0.009841125123581703
